$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the row for "CQRules:CQBP-84--dependencies" (rule retired; row 38)
$ws.Rows.Item(38).Delete()

# 2) Insert a row at 36 to hold the relocated "CloudServiceIncompatibleWorkflowProcess"
#    rule (moved up into the Bug/Critical-ish block, now with severity Blocker)
$ws.Rows.Item(36).Insert()
$ws.Range("A36").Value = "CloudServiceIncompatibleWorkflowProcess"
$ws.Range("B36").Value = "Usage of Cloud Service Incompatible Workflow Processes"
$ws.Range("C36").Value = "Bug"
$ws.Range("D36").Value = "Blocker"
$ws.Range("E36").Value = "aem,cloud-service-compatibility"

# 3) Remove the old duplicate "CloudServiceIncompatibleWorkflowProcess" / Major row
#    (after the insert above it now sits at row 45)
$ws.Rows.Item(45).Delete()

# 4) Insert a new row before "ClientlibProxyResource" (row 70) for the new
#    "IndexDamAssetLucene" rule and populate it
$ws.Rows.Item(70).Insert()
$ws.Range("A70").Value = "IndexDamAssetLucene"
$ws.Range("B70").Value = "Index customizations of the damAssetLucene Oak index should be properly structured."
$ws.Range("C70").Value = "Bug"
$ws.Range("D70").Value = "Minor"
$ws.Range("E70").Value = "aem,cloud-service-compatibility"

# 5) Match the author's final selection/cursor position
$ws.Range("E70").Select()

Write-Output "edit complete"
